$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row for the burden-mode scoring block
$ws.Range("A12").Value = "BURDEN_MODE_Defeat_score="
$ws.Range("B12").Value = "enemy_speed"
$ws.Range("C12").Value = "level"
$ws.Range("D12").Value = "negative effects"
$ws.Range("E12").Value = "result"

# Data row with the new scoring formula
$ws.Range("B13").Value = 500
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 0
$ws.Range("E13").Formula = "=POWER(1.25,D13)*(1000-B13)*C13*C13/5"

# Update the selection / view to match the edited workbook
$ws.Range("D12").Select()
